$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity-log rows to append (A: user, B: activity, C: date, D: time).
# Leading apostrophes force literal-text entry so date/time-looking strings
# (e.g. "02/09/24", "09:40:38") aren't auto-coerced into numeric date/time
# serials, matching the existing inline-string cells in the sheet. The
# Style reset clears the "quote prefix" cell style that entry mode leaves
# behind so the new cells keep the sheet's default (unstyled) formatting.
$data = @(
    @("q", "LOG-IN", "02/09/24", "09:40:38"),
    @("q", "LOG-IN", "02/09/24", "09:41:35"),
    @("q", "LOG-IN", "02/09/24", "09:42:06"),
    @("q", "LOG-IN", "02/09/24", "09:48:25"),
    @("q", "LOG-IN", "02/09/24", "09:55:39"),
    @("q", "LOG-IN", "02/09/24", "09:57:55"),
    @("q", "LOG-IN", "02/09/24", "09:58:31"),
    @("q", "LOG-IN", "02/09/24", "10:00:40"),
    @("q", "LOG-IN", "02/09/24", "10:09:29"),
    @("q", "LOG-IN", "02/09/24", "10:17:21"),
    @("q", "LOG-IN", "02/09/24", "10:19:06"),
    @("q", "LOG-IN", "02/09/24", "10:28:15"),
    @("q", "LOG-IN", "02/09/24", "10:29:17")
)

$startRow = 152
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = "'" + $row[$c - 1]
        $cell.Style = "Normal"
    }
}
